$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 to hold the combined tuple-like string representation
$ws.Range("A2").Value = "('Elemental Shaman', ['Token Creature " + [char]0x2014 + " Elemental Shaman', '3/1'])"

# Remove rows 3 and 4 entirely so the used range shrinks back to A1:A2
$ws.Range("A3:A4").EntireRow.Delete()
